# "new Madigan bike hours" - update the "Riders" (column C) and
# "Average" (column D) values on the Ridership sheet for each day of
# the week with the revised figures. The embedded line chart's series
# (Ridership!$C$2:$C$8 and Ridership!$D$2:$D$8) read directly from
# these worksheet cells, so it reflects the new numbers whenever the
# chart/workbook is recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 235
$ws.Range("D2").Value = 225.59

$ws.Range("C3").Value = 186
$ws.Range("D3").Value = 210.29

$ws.Range("C4").Value = 207
$ws.Range("D4").Value = 212.96

$ws.Range("C5").Value = 283
$ws.Range("D5").Value = 241.22

$ws.Range("C6").Value = 245
$ws.Range("D6").Value = 242.92

$ws.Range("C7").Value = 144
$ws.Range("D7").Value = 114.28

$ws.Range("C8").Value = 85
$ws.Range("D8").Value = 92.73999999999999
